$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.137218475341797
$ws.Range("C2").Value = 5.413793087005615
$ws.Range("D2").Value = 11.997179985046387
$ws.Range("E2").Value = 46.42856979370117
